$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the login test-data rows (2-5) so every row uses the same
# credentials (test@kennect.io / Qwerty@1234) and only the browser column
# differs: chrome, edge, firefox, chrome. This folds the old "test#kennect.io"
# / "Qwerty@12345" variants out of the sheet (and, transitively, out of the
# shared-strings table once they are no longer referenced anywhere).
$ws.Range("A3").Value = "test@kennect.io"
$ws.Range("B3").Value = "Qwerty@1234"
$ws.Range("C3").Value = "edge"

$ws.Range("A4").Value = "test@kennect.io"
$ws.Range("B4").Value = "Qwerty@1234"
$ws.Range("C4").Value = "firefox"

$ws.Range("A5").Value = "test@kennect.io"
$ws.Range("B5").Value = "Qwerty@1234"
$ws.Range("C5").Value = "chrome"

# Hyperlinks: only A2, B2 and A4 should keep a mailto: hyperlink afterwards;
# drop the rest (A3, B3, B4, A5, B5). The host engine only exposes a
# whole-collection Hyperlinks.Delete(), so clear everything and recreate the
# three that must survive.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test@kennect.io") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Qwerty@1234") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:test@kennect.io") | Out-Null

# Hyperlinks.Add() stamps a fresh (duplicate) "Hyperlink" style record;
# reapply the named style so these cells keep using the original one.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"

# Move the active selection to B5.
$ws.Range("B5").Select() | Out-Null
